$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

# --- Insert a new row above the "jogging_km" row (row 7) to host the new
#     "jogging_at_all" multiple-choice question. This pushes jogging_km,
#     jogging_min, mood, drinks, submit and note rows down by one, and the
#     existing (already blank) gap row that used to sit at row 9 slides down
#     to become the new row 10 -- exactly where the new "feelgood" question
#     needs to live. ---
$ws.Rows.Item(7).Insert()

# New row 7: jogging_at_all (mc)
$ws.Cells.Item(7, 3).Value = "mc"
$ws.Cells.Item(7, 4).Value = "jogging_at_all"
$ws.Cells.Item(7, 6).Value = "Did you go jogging today?"
$ws.Cells.Item(7, 7).WrapText = $true
$ws.Cells.Item(7, 7).Value = "yes"
$ws.Cells.Item(7, 8).WrapText = $true
$ws.Cells.Item(7, 8).Value = "no"
$ws.Cells.Item(7, 10).Clear()
$ws.Rows.Item(7).RowHeight = 15

# --- Header row: rename "skipif" column to the new JS-evaluated "showif" ---
$ws.Cells.Item(1, 9).Value = "showif"

# Rows 8 & 9 (jogging_km / jogging_min) now need the new "showif" formula
# text in column I, and a taller row height to match the new content.
$ws.Cells.Item(8, 9).WrapText = $true
$ws.Cells.Item(8, 9).Value = "tail(jogging_at_all, 1) == 1"
$ws.Rows.Item(8).RowHeight = 45
$ws.Rows.Item(9).RowHeight = 45

# Row 10 is presently an empty gap row (it used to be row 9 before the
# insert above); populate it with the new "feelgood" mc_button question.
$ws.Cells.Item(10, 3).WrapText = $true
$ws.Cells.Item(10, 3).Value = "mc_button"
$ws.Cells.Item(10, 4).WrapText = $true
$ws.Cells.Item(10, 4).Value = "feelgood"
$ws.Cells.Item(10, 6).WrapText = $true
$ws.Cells.Item(10, 6).NumberFormat = "@"
$ws.Cells.Item(10, 6).Value = "Did it feel good?"
$ws.Cells.Item(10, 7).WrapText = $true
$ws.Cells.Item(10, 7).Value = "yes"
$ws.Cells.Item(10, 8).WrapText = $true
$ws.Cells.Item(10, 8).Value = "no"
$ws.Cells.Item(10, 9).WrapText = $true
$ws.Cells.Item(10, 9).Value = "tail(jogging_at_all, 1) == 1"
$ws.Rows.Item(10).RowHeight = 45

# Row 9 also gets the showif formula (string already interned above).
$ws.Cells.Item(9, 9).WrapText = $true
$ws.Cells.Item(9, 9).Value = "tail(jogging_at_all, 1) == 1"

# --- Update selection to match the author's final cursor position ---
$ws.Range("I10").Select()
